$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Buenos Aires Innovation Park (City of Buenos Aires Government)"
# which sits at row 4 (A4). Deleting the whole row shifts the rows below it
# up by one and Excel drops the now-unused shared string on save.
$ws.Rows.Item(4).Delete()
